$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.488.04"
$ws.Range("E2").Value = "  +1.41%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.032.49"
$ws.Range("E3").Value = "  +2.40%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "256.69"
$ws.Range("E5").Value = "  +4.64%  "

# Row 6
$ws.Range("E6").Value = "  -0.82%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.37"
$ws.Range("E8").Value = "  -6.05%  "

# Row 9
$ws.Range("E9").Value = "  +1.22%  "

# Row 10
$ws.Range("E10").Value = "  -0.63%  "

# Row 11
$ws.Range("E11").Value = "  -1.80%  "

# Row 12
$ws.Range("E12").Value = "  -1.24%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.333.02"
$ws.Range("E13").Value = "  +2.57%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.822"
$ws.Range("E14").Value = "  -2.87%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.43"
$ws.Range("E15").Value = "  -3.50%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.36"
$ws.Range("E16").Value = "  -2.17%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.045.01"
$ws.Range("E17").Value = "  +3.32%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.447.69"
$ws.Range("E18").Value = "  +1.64%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.00"
$ws.Range("E19").Value = "  -0.44%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0854"
$ws.Range("E20").Value = "  -0.93%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.21"
$ws.Range("E21").Value = "  +0.57%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.34"
$ws.Range("E22").Value = "  -0.24%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.64"
$ws.Range("E23").Value = "  +4.84%  "

# Row 24
$ws.Range("E24").Value = "  -0.03%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  -1.15%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.14"
$ws.Range("E26").Value = "  -1.73%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.48"
$ws.Range("E27").Value = "  +0.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.136"
$ws.Range("E28").Value = "  -7.60%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.90"
$ws.Range("E29").Value = "  +1.73%  "

# Row 30
$ws.Range("E30").Value = "  -1.10%  "

# Row 31
$ws.Range("E31").Value = "  -1.00%  "

# Row 32
$ws.Range("E32").Value = "  +7.32%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.73"
$ws.Range("E33").Value = "  -3.32%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.55"
$ws.Range("E34").Value = "  +0.03%  "

# Row 35
$ws.Range("E35").Value = "  +7.91%  "

# Row 36
$ws.Range("E36").Value = "  -0.01%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.41"
$ws.Range("E37").Value = "  +1.45%  "

# Row 38
$ws.Range("E38").Value = "  +2.28%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.37"
$ws.Range("E39").Value = "  -3.45%  "

# Row 41
$ws.Range("E41").Value = "  -2.90%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0216"
$ws.Range("E42").Value = "  +1.28%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.19"
$ws.Range("E43").Value = "  +0.76%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.404.31"
$ws.Range("E44").Value = "  +2.27%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.13"
$ws.Range("E45").Value = "  -2.55%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.84"
$ws.Range("E46").Value = "  +0.59%  "

# Row 47
$ws.Range("E47").Value = "  +0.88%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.35"
$ws.Range("E48").Value = "  +1.17%  "

# Row 49
$ws.Range("E49").Value = "  +1.69%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.01"
$ws.Range("E50").Value = "  +1.15%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.223.67"
$ws.Range("E51").Value = "  +2.55%  "
